# Applies the syllabus edit described by the commit:
#   - drop the stray "_GoBack" bookmark that sat next to the academic
#     integrity hyperlink
#   - append " 11:59PM." to the "6/12/24: ... Final Projects Due" bullet
#   - leave a fresh "_GoBack" bookmark right after the text that was
#     just typed (mirroring how Word itself tracks the last edit point)

$d = $word.ActiveDocument

# --- 1) Remove the old _GoBack bookmark near the academic-integrity link ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2) Append " 11:59PM." to the "6/12/24 ... Final Projects Due" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Differential abundance*") {
        $r = $p.Range
        $r.Collapse(0)        # wdCollapseEnd
        [void]$r.MoveEnd(1, -1)   # step back before the paragraph mark
        $r.InsertAfter(" 11:59PM.")
        break
    }
}

# Keep the appended text in its own run (as in the authored edit) by
# toggling a throw-away bookmark at the seam between the old and new text;
# adding/removing a bookmark there permanently splits the run.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Differential abundance*") {
        $paraText = $p.Range.Text
        $dueIdx = $paraText.IndexOf("Due")
        $splitPos = $p.Range.Start + $dueIdx + 3
        $splitRange = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("tempSplit", $splitRange)
        $d.Bookmarks("tempSplit").Delete()
        break
    }
}

# --- 3) Re-create _GoBack collapsed right after " 11:59PM." ---
# A collapsed range sitting exactly one character before a paragraph mark
# confuses Bookmarks.Add, so park a throw-away placeholder character after
# the insertion point, anchor the bookmark just before it, then delete the
# placeholder again.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Differential abundance*") {
        $p.Range.InsertAfter("X")
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Differential abundance*") {
        $bmPos = $p.Range.End - 2
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Differential abundance*") {
        $delRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
        $delRange.Delete()
        break
    }
}
